$d = $word.ActiveDocument

# 1) Replace equipment model text
$d.Content.Find.Execute("MEDISONIC MODELO H60 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MINDRAY MODELO DC – N3 ", 2)

# 2) Fix "ANEXO IZQUIERDO" paragraph: merge "ANEXO IZQUIERDO" + ":" into one bold run,
#    and merge "." + " Trompa libre..." into one italic run.
$d.Content.Find.Execute("ANEXO IZQUIERDO:.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ANEXO IZQUIERDO:.", 2)
